# Update the inline sample data: replace the old single sample row with the
# new calibration rows (157/158/159) and their d18O readings, extending the
# used range from A1:B2 to A1:B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sample IDs in column A must stay text (they look numeric), so force the
# cells to a text format before writing the values - mirrors how Excel
# keeps IDs like "157" from being auto-converted to numbers.
$ws.Range("A2:A4").NumberFormat = "@"

# Row 2: old placeholder sample -> sample 157
$ws.Range("A2").Value = "157"
$ws.Range("B2").Value = 31.1

# Row 3 (new): sample 158
$ws.Range("A3").Value = "158"
$ws.Range("B3").Value = 30.2

# Row 4 (new): sample 159
$ws.Range("A4").Value = "159"
$ws.Range("B4").Value = 29.8
